$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- 1) Revise the published GDP-per-Capita series for the years already on
#     the sheet (rows 2-192, column E). The source values were restated
#     (e.g. 1820: 570 -> 909 ... 2010: 9372.20940704 -> 13343.546749821);
#     only the rows whose figure actually changed are listed below.
#     Column E stores these as text (matching the original workbook), so the
#     range is forced to Text format before the assignment and then reset to
#     the sheet's default (unstyled) look.
$revisedValues = @(
    @(2, "909"),
    @(52, "969"),
    @(72, "1250"),
    @(95, "1341"),
    @(111, "1264"),
    @(120, "1317"),
    @(132, "1302"),
    @(133, "1353"),
    @(134, "1385"),
    @(135, "1490"),
    @(136, "1431"),
    @(137, "1506"),
    @(138, "1482"),
    @(139, "1451"),
    @(140, "1457"),
    @(141, "1581"),
    @(142, "1718"),
    @(143, "1753"),
    @(144, "1831"),
    @(145, "1921"),
    @(146, "1991"),
    @(147, "2085"),
    @(148, "2251"),
    @(149, "2369"),
    @(150, "2488"),
    @(151, "2608"),
    @(152, "2700"),
    @(153, "2750"),
    @(154, "2786"),
    @(155, "2987"),
    @(156, "3045"),
    @(157, "3123"),
    @(158, "3333"),
    @(159, "3585"),
    @(160, "3861"),
    @(161, "3979"),
    @(162, "4071"),
    @(163, "4230"),
    @(164, "4375"),
    @(165, "4538"),
    @(166, "4720"),
    @(167, "4862"),
    @(168, "5053"),
    @(169, "5453"),
    @(170, "6089"),
    @(171, "6735"),
    @(172, "7385"),
    @(173, "7849.74481802388"),
    @(174, "8409.0209422082"),
    @(175, "8964.29215998945"),
    @(176, "9499.51752452374"),
    @(177, "10080.2082631763"),
    @(178, "10458.3586098259"),
    @(179, "9995.61821964787"),
    @(180, "9077.7477967745"),
    @(181, "9347.0609300166"),
    @(182, "9627.10431781556"),
    @(183, "9818.72053622492"),
    @(184, "10276.398308003"),
    @(185, "10868.8355383604"),
    @(186, "11408.3257048843"),
    @(187, "11742.1280936609"),
    @(188, "12180.7113589908"),
    @(189, "12694.7146347803"),
    @(190, "12768.1290769284"),
    @(191, "12543.3986128863"),
    @(192, "13343.546749821")
)

foreach ($entry in $revisedValues) {
    $row = $entry[0]
    $value = $entry[1]
    $cell = $ws.Range("E$row")
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# --- 2) Append the newly published years (2011-2016) as new rows.
$newYears = @(
    @(193, 2011, "13318"),
    @(194, 2012, "14227"),
    @(195, 2013, "14562"),
    @(196, 2014, "14642"),
    @(197, 2015, "15020"),
    @(198, 2016, "15454")
)

foreach ($entry in $newYears) {
    $row = $entry[0]
    $year = $entry[1]
    $value = $entry[2]

    $ws.Range("A$row").Value = 764
    $ws.Range("B$row").Value = "Thailand"
    $ws.Range("C$row").Value = "GDP per Capita"
    $ws.Range("D$row").Value = $year

    $cell = $ws.Range("E$row")
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}
